# wordform parsing behaviour changed django_toolbar added
#
# The only real content change in this commit is the rich-text string in
# C4: it drops the "[Hev] [Ro]" tag from the first run and the trailing
# " [Ro]" from the second run, going from
#   "aBi | api [Hev] [Ro] | app'" + "ee [Ro]"
# to
#   "aBi | api | app'" + "ee"
# (everything else in the sharedStrings table just shifts index because
# this entry is re-inserted at the end of the table).  The sheet also
# records that C7 was the active cell when the file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c4 = $ws.Range("C4")

# Replace the full text, then re-apply the distinct formatting run ("ee")
# that the second <r> in the original rich string carried, so the cell
# keeps two runs instead of collapsing to a single plain string.
$c4.Value = "aBi | api | app'ee"

$tail = $c4.Characters(17, 2)          # "ee" - the old second run's text
$tail.Font.Name = "Calibri"
$tail.Font.Size = 11
$tail.Font.ColorIndex = -4105          # xlColorIndexAutomatic - no explicit color, as in the source run

# Record the selection that was active when the workbook was saved.
$ws.Range("C7").Select() | Out-Null
